$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $row = [int]($cellRef -replace "[A-Z]+", "")
    $styleSrc = $ws.Range("B" + $row)
    $target = $ws.Range($cellRef)
    $target.Value = "'" + $text
    $target.Style = $styleSrc.Style
}

Set-TextValue "D2" "42.542.13"
Set-TextValue "E2" "  +1.90%  "
Set-TextValue "D3" "2.290.30"
Set-TextValue "E3" "  +0.91%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "E5" "  +1.38%  "
Set-TextValue "D6" "97.48"
Set-TextValue "E6" "  +5.91%  "
Set-TextValue "D7" "0.532"
Set-TextValue "E7" "  +0.54%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "D9" "0.497"
Set-TextValue "E9" "  +3.47%  "
Set-TextValue "D10" "36.29"
Set-TextValue "E10" "  +12.33%  "
Set-TextValue "E11" "  +1.14%  "
Set-TextValue "E12" "  -1.49%  "
Set-TextValue "D13" "6.74"
Set-TextValue "E13" "  +2.18%  "
Set-TextValue "D14" "2.649.36"
Set-TextValue "E14" "  +1.06%  "
Set-TextValue "D15" "14.60"
Set-TextValue "E15" "  +2.85%  "
Set-TextValue "D16" "2.304.83"
Set-TextValue "E16" "  +0.92%  "
Set-TextValue "D17" "0.806"
Set-TextValue "E17" "  +5.68%  "
Set-TextValue "D18" "42.467.97"
Set-TextValue "E18" "  +1.88%  "
Set-TextValue "D19" "12.75"
Set-TextValue "E19" "  +1.48%  "
Set-TextValue "D20" "0.0₃0920"
Set-TextValue "E20" "  +1.99%  "
Set-TextValue "D21" "6.03"
Set-TextValue "E21" "  +2.05%  "
Set-TextValue "D22" "67.99"
Set-TextValue "E22" "  +1.74%  "
Set-TextValue "D23" "243.37"
Set-TextValue "E23" "  +1.33%  "
Set-TextValue "D24" "2.61"
Set-TextValue "E24" "  +1.11%  "
Set-TextValue "E25" "  +2.63%  "
Set-TextValue "E26" "  -0.16%  "
Set-TextValue "D27" "24.02"
Set-TextValue "E27" "  -0.12%  "
Set-TextValue "D28" "37.49"
Set-TextValue "E28" "  +9.83%  "
Set-TextValue "E29" "  +0.96%  "
Set-TextValue "D30" "2.12"
Set-TextValue "E30" "  +2.44%  "
Set-TextValue "D31" "161.17"
Set-TextValue "E31" "  +0.37%  "
Set-TextValue "D32" "5.32"
Set-TextValue "E32" "  +1.65%  "
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  +0.00%  "
Set-TextValue "D34" "3.14"
Set-TextValue "E34" "  +4.41%  "
Set-TextValue "E35" "  +1.59%  "
Set-TextValue "D36" "17.38"
Set-TextValue "E36" "  +3.06%  "
Set-TextValue "E37" "  +3.46%  "
Set-TextValue "D38" "1.88"
Set-TextValue "E38" "  +5.00%  "
Set-TextValue "E39" "  +0.20%  "
Set-TextValue "E40" "  -0.14%  "
Set-TextValue "D41" "4.21"
Set-TextValue "E41" "  +6.81%  "
Set-TextValue "E42" "  +17.26%  "
Set-TextValue "D43" "2.004.69"
Set-TextValue "E43" "  -1.94%  "
Set-TextValue "B44" "VeChain"
Set-TextValue "C44" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D44" "0.0288"
Set-TextValue "E44" "  +3.56%  "
Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "19.23"
Set-TextValue "E45" "  -0.10%  "
Set-TextValue "D46" "3.03"
Set-TextValue "E46" "  +5.23%  "
Set-TextValue "D47" "10.24"
Set-TextValue "E47" "  -1.25%  "
Set-TextValue "D48" "53.84"
Set-TextValue "E48" "  +4.05%  "
Set-TextValue "E49" "  +0.89%  "
Set-TextValue "D50" "72.91"
Set-TextValue "E50" "  +0.40%  "
Set-TextValue "E51" "  -0.51%  "
